$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "9802.IB" (row 33), shifting the remaining cells up.
$found = $ws.Range("A33")
$found.EntireRow.Delete()

# Match the resulting selection left behind in Excel after the delete
# (Excel scrolled the view so row 19 is the top visible row; this COM
# runtime does not persist ScrollRow/topLeftCell to the saved file when
# no panes are frozen, so only the selection itself is reproduced here).
$ws.Range("A33:XFD33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
